$wb = $excel.ActiveWorkbook

# --- Sheet "bug" (sheet1) ---
$ws1 = $wb.Worksheets.Item("bug")

$ws1.Range("C47").Value = "N/A"
$ws1.Range("D47").Value = "android"
$ws1.Range("E47").Value = 20151215
$ws1.Range("F47").Value = "XOPlayer 部分mp4用ffextrator打开很慢(system 正常)"
$ws1.Range("G47").Value = "fixed"
$ws1.Range("H47").Value = "获取 sps pps错误(AVCC解析)"

$ws1.Range("C48").Value = "N/A"
$ws1.Range("D48").Value = "android"
$ws1.Range("E48").Value = 20151216
$ws1.Range("F48").Value = "XOPlayer 本地视频拖动慢"
$ws1.Range("G48").Value = "tracking"

$ws1.Range("C49").Value = "N/A"
$ws1.Range("D49").Value = "android"
$ws1.Range("E49").Value = 20151216
$ws1.Range("F49").Value = "XOPlayer PPTV视频本地播放SystemMediaExtractor声音不对"
$ws1.Range("G49").Value = "tracking"

$ws1.Range("C50").Value = "N/A"
$ws1.Range("D50").Value = "android"
$ws1.Range("E50").Value = 20151216
$ws1.Range("F50").Value = "XOPlayer PPTV视频播放不流畅"
$ws1.Range("G50").Value = "tracking"

$ws1.Activate()
$excel.ActiveWindow.ScrollRow = 31
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws1.Range("F50").Select()

# --- Sheet "newfeature" (sheet2) ---
$ws2 = $wb.Worksheets.Item("newfeature")

$ws2.Range("B39").Value = 38
$ws2.Range("C39").Value = "meetsdk"
$ws2.Range("D39").Value = "android"
$ws2.Range("E39").Value = "FFPlayer FFExtractor分离"
$ws2.Range("F39").Value = "done"

$ws2.Range("B40").Value = 39
$ws2.Range("C40").Value = "meetsdk"
$ws2.Range("D40").Value = "all"
$ws2.Range("E40").Value = "ffmpeg 裁剪"
$ws2.Range("F40").Value = "TBD"

$ws2.Range("B41").Value = 40

$ws2.Activate()
$excel.ActiveWindow.ScrollRow = 19
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws2.Range("C41").Select()

$ws1.Activate()
